# Generate Report for Handoff
# Updates the "b.md" rows across the Overview / zh-cn / de-de sheets to
# reflect that a new handoff (b.*.xlf) has been generated, while the
# handback file on record is still the stale one -> error detail note.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84bf0ddc1217c010f496ae98a4bf72ba342732bf/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/82c1eb107932c1af340ac365181816b23de151ae/e2e/b.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-26 07:30:51"

# ---- zh-cn sheet (row 3 = b.md) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-10-26 07:30:38"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# ---- de-de sheet (row 3 = b.md) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-10-26 07:30:51"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 40
